$d = $word.ActiveDocument
$d.Content.Find.Execute("Miloš Jovanović", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Miloš Ćirković", 2)
